# "Список дел" (To-do list) update:
#  - mark 3 tasks (row 4 "ctrl+c ...", row 7 "добавить сохранения и загрузку...",
#    row 17 "Пофиксить краш сервера ...") as done ("Cделано", green fill)
#  - row 4 also gets a "В процессе" (in progress, gold fill) marker in column D
#  - a brand new task row (18) is appended: "Добавить новый тип предметов (Документы)"
#    with the usual "Не сделано" (red) status
#  - view resets: selection moves to B4, the old scroll-anchor (topLeftCell=A9) goes away
#  - column D gets an explicit width so the new "В процессе" tag is readable

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlThemeColorAccent6 (green) == OOXML theme index 9 -> "Cделано"
$themeAccent6 = 10
# xlThemeColorAccent4 (gold) == OOXML theme index 7 -> "В процессе"
$themeAccent4 = 8

# --- Row 4: "ctrl+c пизда серверу если клиент закрыть" -> Cделано, + В процессе in D4 ---
$ws.Range("B4").Value = "Cделано"
$ws.Range("B4").Interior.ThemeColor = $themeAccent6

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (copy border/alignment template)
$ws.Range("D4").Value = "В процессе"
$ws.Range("D4").Interior.ThemeColor = $themeAccent4

# --- Row 7: "добавить сохранения и загрузку(Бинарник)" -> Cделано ---
$ws.Range("B7").Value = "Cделано"
$ws.Range("B7").Interior.ThemeColor = $themeAccent6

# --- Row 17: "Пофиксить краш сервера при попытки зайти с одинаковым ником" -> Cделано ---
$ws.Range("B17").Value = "Cделано"
$ws.Range("B17").Interior.ThemeColor = $themeAccent6

# --- New row 18: new task, still "Не сделано" ---
$ws.Range("A8").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "Добавить новый тип предметов (Документы)"

$ws.Range("B8").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "Не сделано"

# --- column D needs an explicit width now that it carries visible text ---
$ws.Columns.Item(4).ColumnWidth = 10

# --- view: selection moves to B4, scroll anchor resets to default ---
$ws.Range("B4").Select() | Out-Null
